$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "28.051.92"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.872.99"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.93%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.33%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "312.86"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("E6").Value = "  +0.27%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.5116"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +1.35%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3824"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.81%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.08274"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -10.49%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.114"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -1.19%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "41.66"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.40%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "6.228"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.882.08"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.87%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "20.51"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.38%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.212"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -1.04%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.00001095"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "91.01"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.34%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06649"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "17.98"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  +0.20%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.056"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.48%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "28.096.46"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  -1.97%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.268"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -2.17%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.595"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +2.13%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.097.74"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.82%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "157.50"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.54%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "20.60"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.97%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "125.85"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.74%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.1057"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("E32").Value = "  -2.89%  "
$ws.Range("E33").Value = "  -0.02%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.608"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.08%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "9.674"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +2.10%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.02451"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.95%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.06574"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.45%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.2168"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.27%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.214"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.6487"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.248"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -7.23%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "11.35"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.50%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "4.887"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -1.58%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.6137"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +1.60%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "13.12"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.25%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "1.297"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.24%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "3.674"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.38%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.012"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.55%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.218"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +2.02%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "120.98"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.00%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "80.60"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +2.16%  "
